$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before D; existing D:K data shifts to E:L
$ws.Columns("D").Insert()

# Copy cell formatting (number format, font, alignment) from column E (the shifted
# original column D) into the newly inserted column D, for every row that has data
$ws.Range("E7:E35,E38:E77,E80:E102").Copy()
$ws.Range("D7:D35,D38:D77,D80:D102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Write the refreshed quarterly figures (columns D through K); column L already holds
# the correct value because it was shifted automatically from the old column K above.
$row7 = New-Object 'object[,]' 1,8
$row7[0,0] = 43373
$row7[0,1] = 43281
$row7[0,2] = 43190
$row7[0,3] = 43100
$row7[0,4] = 43008
$row7[0,5] = 42916
$row7[0,6] = 42825
$row7[0,7] = 42735
$ws.Range("D7:K7").Value2 = $row7
$row8 = New-Object 'object[,]' 1,8
$row8[0,0] = 8000
$row8[0,1] = 11600
$row8[0,2] = 5700
$row8[0,3] = 6800
$row8[0,4] = 8300
$row8[0,5] = 10900
$row8[0,6] = 5000
$row8[0,7] = 7400
$ws.Range("D8:K8").Value2 = $row8
$row9 = New-Object 'object[,]' 1,8
$row9[0,0] = 5400
$row9[0,1] = 7100
$row9[0,2] = 3100
$row9[0,3] = 4100
$row9[0,4] = 4800
$row9[0,5] = 6900
$row9[0,6] = 3300
$row9[0,7] = 5900
$ws.Range("D9:K9").Value2 = $row9
$row10 = New-Object 'object[,]' 1,8
$row10[0,0] = 2600
$row10[0,1] = 4500
$row10[0,2] = 2600
$row10[0,3] = 2700
$row10[0,4] = 3500
$row10[0,5] = 4000
$row10[0,6] = 1700
$row10[0,7] = 1500
$ws.Range("D10:K10").Value2 = $row10
$row12 = New-Object 'object[,]' 1,8
$row12[0,0] = 1900
$row12[0,1] = 3500
$row12[0,2] = 1800
$row12[0,3] = 1800
$row12[0,4] = 1300
$row12[0,5] = 3200
$row12[0,6] = 1700
$row12[0,7] = 2200
$ws.Range("D12:K12").Value2 = $row12
$row13 = New-Object 'object[,]' 1,8
$row13[0,0] = 0
$row13[0,1] = 0
$row13[0,2] = 0
$row13[0,3] = 0
$row13[0,4] = 0
$row13[0,5] = 0
$row13[0,6] = 0
$row13[0,7] = 0
$ws.Range("D13:K13").Value2 = $row13
$row14 = New-Object 'object[,]' 1,8
$row14[0,0] = 0
$row14[0,1] = 0
$row14[0,2] = 0
$row14[0,3] = 0
$row14[0,4] = 0
$row14[0,5] = 0
$row14[0,6] = 0
$row14[0,7] = 0
$ws.Range("D14:K14").Value2 = $row14
$row15 = New-Object 'object[,]' 1,8
$row15[0,0] = 0
$row15[0,1] = 0
$row15[0,2] = 0
$row15[0,3] = 0
$row15[0,4] = 0
$row15[0,5] = 0
$row15[0,6] = 0
$row15[0,7] = 0
$ws.Range("D15:K15").Value2 = $row15
$row17 = New-Object 'object[,]' 1,8
$row17[0,0] = 11100
$row17[0,1] = 17000
$row17[0,2] = 8200
$row17[0,3] = 9800
$row17[0,4] = 9400
$row17[0,5] = 16700
$row17[0,6] = 7800
$row17[0,7] = 10800
$ws.Range("D17:K17").Value2 = $row17
$row18 = New-Object 'object[,]' 1,8
$row18[0,0] = -3100
$row18[0,1] = -5400
$row18[0,2] = -2500
$row18[0,3] = -3000
$row18[0,4] = -1100
$row18[0,5] = -5800
$row18[0,6] = -2700
$row18[0,7] = -3400
$ws.Range("D18:K18").Value2 = $row18
$row20 = New-Object 'object[,]' 1,8
$row20[0,0] = -1100
$row20[0,1] = 200
$row20[0,2] = 700
$row20[0,3] = 300
$row20[0,4] = -100
$row20[0,5] = -100
$row20[0,6] = 0
$row20[0,7] = -100
$ws.Range("D20:K20").Value2 = $row20
$row21 = New-Object 'object[,]' 1,8
$row21[0,0] = -3300
$row21[0,1] = -3300
$row21[0,2] = -800
$row21[0,3] = -1700
$row21[0,4] = -200
$row21[0,5] = -4200
$row21[0,6] = -2000
$row21[0,7] = -2700
$ws.Range("D21:K21").Value2 = $row21
$row22 = New-Object 'object[,]' 1,8
$row22[0,0] = 0
$row22[0,1] = 0
$row22[0,2] = 0
$row22[0,3] = 0
$row22[0,4] = 0
$row22[0,5] = 0
$row22[0,6] = 0
$row22[0,7] = 0
$ws.Range("D22:K22").Value2 = $row22
$row23 = New-Object 'object[,]' 1,8
$row23[0,0] = -4300
$row23[0,1] = -5200
$row23[0,2] = -1800
$row23[0,3] = -2700
$row23[0,4] = -1200
$row23[0,5] = -5900
$row23[0,6] = -2700
$row23[0,7] = -3500
$ws.Range("D23:K23").Value2 = $row23
$row24 = New-Object 'object[,]' 1,8
$row24[0,0] = 0
$row24[0,1] = 0
$row24[0,2] = 0
$row24[0,3] = 100
$row24[0,4] = 0
$row24[0,5] = 0
$row24[0,6] = "NA"
$row24[0,7] = 0
$ws.Range("D24:K24").Value2 = $row24
$row25 = New-Object 'object[,]' 1,8
$row25[0,0] = 0
$row25[0,1] = 0
$row25[0,2] = 0
$row25[0,3] = 0
$row25[0,4] = 0
$row25[0,5] = 0
$row25[0,6] = 0
$row25[0,7] = 0
$ws.Range("D25:K25").Value2 = $row25
$row26 = New-Object 'object[,]' 1,8
$row26[0,0] = -4300
$row26[0,1] = -5200
$row26[0,2] = -1800
$row26[0,3] = -2800
$row26[0,4] = -1200
$row26[0,5] = -5900
$row26[0,6] = -2700
$row26[0,7] = -3500
$ws.Range("D26:K26").Value2 = $row26
$row27 = New-Object 'object[,]' 1,8
$row27[0,0] = -4300
$row27[0,1] = -5200
$row27[0,2] = -1800
$row27[0,3] = -2800
$row27[0,4] = -1200
$row27[0,5] = -5900
$row27[0,6] = -2700
$row27[0,7] = -3500
$ws.Range("D27:K27").Value2 = $row27
$row28 = New-Object 'object[,]' 1,8
$row28[0,0] = 0
$row28[0,1] = 0
$row28[0,2] = 0
$row28[0,3] = 0
$row28[0,4] = 0
$row28[0,5] = 0
$row28[0,6] = 0
$row28[0,7] = 0
$ws.Range("D28:K28").Value2 = $row28
$row29 = New-Object 'object[,]' 1,8
$row29[0,0] = 0
$row29[0,1] = 0
$row29[0,2] = 0
$row29[0,3] = 0
$row29[0,4] = 0
$row29[0,5] = 0
$row29[0,6] = 0
$row29[0,7] = 0
$ws.Range("D29:K29").Value2 = $row29
$row30 = New-Object 'object[,]' 1,8
$row30[0,0] = 0
$row30[0,1] = 0
$row30[0,2] = 0
$row30[0,3] = 0
$row30[0,4] = 0
$row30[0,5] = 0
$row30[0,6] = 0
$row30[0,7] = 0
$ws.Range("D30:K30").Value2 = $row30
$row31 = New-Object 'object[,]' 1,8
$row31[0,0] = 0
$row31[0,1] = 0
$row31[0,2] = 0
$row31[0,3] = 0
$row31[0,4] = 0
$row31[0,5] = 0
$row31[0,6] = 0
$row31[0,7] = 0
$ws.Range("D31:K31").Value2 = $row31
$row32 = New-Object 'object[,]' 1,8
$row32[0,0] = 1100
$row32[0,1] = -200
$row32[0,2] = -700
$row32[0,3] = -300
$row32[0,4] = 100
$row32[0,5] = 100
$row32[0,6] = 0
$row32[0,7] = 100
$ws.Range("D32:K32").Value2 = $row32
$row33 = New-Object 'object[,]' 1,8
$row33[0,0] = -4300
$row33[0,1] = -5200
$row33[0,2] = -1800
$row33[0,3] = -2800
$row33[0,4] = -1200
$row33[0,5] = -5900
$row33[0,6] = -2700
$row33[0,7] = -3500
$ws.Range("D33:K33").Value2 = $row33
$row34 = New-Object 'object[,]' 1,8
$row34[0,0] = 0
$row34[0,1] = 0
$row34[0,2] = 0
$row34[0,3] = 0
$row34[0,4] = 0
$row34[0,5] = 0
$row34[0,6] = 0
$row34[0,7] = 0
$ws.Range("D34:K34").Value2 = $row34
$row35 = New-Object 'object[,]' 1,8
$row35[0,0] = -4300
$row35[0,1] = -5200
$row35[0,2] = -1800
$row35[0,3] = -2800
$row35[0,4] = -1200
$row35[0,5] = -5900
$row35[0,6] = -2700
$row35[0,7] = -3500
$ws.Range("D35:K35").Value2 = $row35
$row38 = New-Object 'object[,]' 1,8
$row38[0,0] = 43373
$row38[0,1] = 43281
$row38[0,2] = 43190
$row38[0,3] = 43100
$row38[0,4] = 43008
$row38[0,5] = 42916
$row38[0,6] = 42825
$row38[0,7] = 42735
$ws.Range("D38:K38").Value2 = $row38
$row41 = New-Object 'object[,]' 1,8
$row41[0,0] = 3500
$row41[0,1] = 4100
$row41[0,2] = 3500
$row41[0,3] = 8500
$row41[0,4] = 3100
$row41[0,5] = 4500
$row41[0,6] = 7800
$row41[0,7] = 9200
$ws.Range("D41:K41").Value2 = $row41
$row42 = New-Object 'object[,]' 1,8
$row42[0,0] = 11100
$row42[0,1] = 14800
$row42[0,2] = 19800
$row42[0,3] = 15800
$row42[0,4] = 12100
$row42[0,5] = 13100
$row42[0,6] = 13100
$row42[0,7] = 14800
$ws.Range("D42:K42").Value2 = $row42
$row43 = New-Object 'object[,]' 1,8
$row43[0,0] = 5300
$row43[0,1] = 4700
$row43[0,2] = 5400
$row43[0,3] = 5700
$row43[0,4] = 5800
$row43[0,5] = 4500
$row43[0,6] = 4000
$row43[0,7] = 4900
$ws.Range("D43:K43").Value2 = $row43
$row44 = New-Object 'object[,]' 1,8
$row44[0,0] = 12000
$row44[0,1] = 14700
$row44[0,2] = 12700
$row44[0,3] = 10400
$row44[0,4] = 10500
$row44[0,5] = 10700
$row44[0,6] = 10600
$row44[0,7] = 13200
$ws.Range("D44:K44").Value2 = $row44
$row45 = New-Object 'object[,]' 1,8
$row45[0,0] = 1800
$row45[0,1] = 1800
$row45[0,2] = 1600
$row45[0,3] = 1700
$row45[0,4] = 2000
$row45[0,5] = 1900
$row45[0,6] = 2200
$row45[0,7] = 2000
$ws.Range("D45:K45").Value2 = $row45
$row46 = New-Object 'object[,]' 1,8
$row46[0,0] = 33800
$row46[0,1] = 40100
$row46[0,2] = 43000
$row46[0,3] = 42100
$row46[0,4] = 33500
$row46[0,5] = 34700
$row46[0,6] = 37700
$row46[0,7] = 44000
$ws.Range("D46:K46").Value2 = $row46
$row47 = New-Object 'object[,]' 1,8
$row47[0,0] = 300
$row47[0,1] = 1200
$row47[0,2] = 1500
$row47[0,3] = 400
$row47[0,4] = 300
$row47[0,5] = 300
$row47[0,6] = 200
$row47[0,7] = 200
$ws.Range("D47:K47").Value2 = $row47
$row48 = New-Object 'object[,]' 1,8
$row48[0,0] = 31300
$row48[0,1] = 29800
$row48[0,2] = 30100
$row48[0,3] = 33900
$row48[0,4] = 31000
$row48[0,5] = 30300
$row48[0,6] = 30200
$row48[0,7] = 27600
$ws.Range("D48:K48").Value2 = $row48
$row49 = New-Object 'object[,]' 1,8
$row49[0,0] = 1600
$row49[0,1] = 1500
$row49[0,2] = 1300
$row49[0,3] = 1200
$row49[0,4] = 1200
$row49[0,5] = 1100
$row49[0,6] = 1000
$row49[0,7] = 1000
$ws.Range("D49:K49").Value2 = $row49
$row50 = New-Object 'object[,]' 1,8
$row50[0,0] = 0
$row50[0,1] = 0
$row50[0,2] = 0
$row50[0,3] = 0
$row50[0,4] = 0
$row50[0,5] = 0
$row50[0,6] = 0
$row50[0,7] = 0
$ws.Range("D50:K50").Value2 = $row50
$row51 = New-Object 'object[,]' 1,8
$row51[0,0] = 0
$row51[0,1] = 0
$row51[0,2] = 0
$row51[0,3] = 0
$row51[0,4] = 0
$row51[0,5] = 0
$row51[0,6] = 0
$row51[0,7] = 0
$ws.Range("D51:K51").Value2 = $row51
$row52 = New-Object 'object[,]' 1,8
$row52[0,0] = 100
$row52[0,1] = 100
$row52[0,2] = 100
$row52[0,3] = 100
$row52[0,4] = 100
$row52[0,5] = 100
$row52[0,6] = 100
$row52[0,7] = 100
$ws.Range("D52:K52").Value2 = $row52
$row53 = New-Object 'object[,]' 1,8
$row53[0,0] = 0
$row53[0,1] = 0
$row53[0,2] = 0
$row53[0,3] = 0
$row53[0,4] = 0
$row53[0,5] = 0
$row53[0,6] = 0
$row53[0,7] = 0
$ws.Range("D53:K53").Value2 = $row53
$row54 = New-Object 'object[,]' 1,8
$row54[0,0] = 67100
$row54[0,1] = 72700
$row54[0,2] = 76000
$row54[0,3] = 75200
$row54[0,4] = 66000
$row54[0,5] = 66500
$row54[0,6] = 69100
$row54[0,7] = 72900
$ws.Range("D54:K54").Value2 = $row54
$row57 = New-Object 'object[,]' 1,8
$row57[0,0] = 3600
$row57[0,1] = 3700
$row57[0,2] = 3200
$row57[0,3] = 3400
$row57[0,4] = 2900
$row57[0,5] = 2500
$row57[0,6] = 3000
$row57[0,7] = 2100
$ws.Range("D57:K57").Value2 = $row57
$row58 = New-Object 'object[,]' 1,8
$row58[0,0] = 1100
$row58[0,1] = 1100
$row58[0,2] = 1200
$row58[0,3] = 1300
$row58[0,4] = 1400
$row58[0,5] = 1400
$row58[0,6] = 1400
$row58[0,7] = 1500
$ws.Range("D58:K58").Value2 = $row58
$row59 = New-Object 'object[,]' 1,8
$row59[0,0] = 3800
$row59[0,1] = 4900
$row59[0,2] = 5700
$row59[0,3] = 2600
$row59[0,4] = 2400
$row59[0,5] = 3200
$row59[0,6] = 3400
$row59[0,7] = 2900
$ws.Range("D59:K59").Value2 = $row59
$row60 = New-Object 'object[,]' 1,8
$row60[0,0] = 8500
$row60[0,1] = 9600
$row60[0,2] = 10000
$row60[0,3] = 7400
$row60[0,4] = 6700
$row60[0,5] = 7000
$row60[0,6] = 7800
$row60[0,7] = 6500
$ws.Range("D60:K60").Value2 = $row60
$row61 = New-Object 'object[,]' 1,8
$row61[0,0] = 18300
$row61[0,1] = 18300
$row61[0,2] = 18400
$row61[0,3] = 18400
$row61[0,4] = 7500
$row61[0,5] = 6700
$row61[0,6] = 6000
$row61[0,7] = 5700
$ws.Range("D61:K61").Value2 = $row61
$row62 = New-Object 'object[,]' 1,8
$row62[0,0] = 300
$row62[0,1] = 300
$row62[0,2] = 200
$row62[0,3] = 100
$row62[0,4] = 100
$row62[0,5] = 200
$row62[0,6] = 200
$row62[0,7] = 300
$ws.Range("D62:K62").Value2 = $row62
$row63 = New-Object 'object[,]' 1,8
$row63[0,0] = 0
$row63[0,1] = 0
$row63[0,2] = 0
$row63[0,3] = 0
$row63[0,4] = 0
$row63[0,5] = 0
$row63[0,6] = 0
$row63[0,7] = 0
$ws.Range("D63:K63").Value2 = $row63
$row64 = New-Object 'object[,]' 1,8
$row64[0,0] = 0
$row64[0,1] = 0
$row64[0,2] = 0
$row64[0,3] = 0
$row64[0,4] = 0
$row64[0,5] = 0
$row64[0,6] = 0
$row64[0,7] = 0
$ws.Range("D64:K64").Value2 = $row64
$row65 = New-Object 'object[,]' 1,8
$row65[0,0] = 0
$row65[0,1] = 0
$row65[0,2] = 0
$row65[0,3] = 0
$row65[0,4] = 0
$row65[0,5] = 0
$row65[0,6] = 0
$row65[0,7] = 0
$ws.Range("D65:K65").Value2 = $row65
$row66 = New-Object 'object[,]' 1,8
$row66[0,0] = 27100
$row66[0,1] = 28300
$row66[0,2] = 28700
$row66[0,3] = 26000
$row66[0,4] = 14300
$row66[0,5] = 13900
$row66[0,6] = 14100
$row66[0,7] = 12500
$ws.Range("D66:K66").Value2 = $row66
$row68 = New-Object 'object[,]' 1,8
$row68[0,0] = 0
$row68[0,1] = 0
$row68[0,2] = 0
$row68[0,3] = 0
$row68[0,4] = 0
$row68[0,5] = 0
$row68[0,6] = 0
$row68[0,7] = 0
$ws.Range("D68:K68").Value2 = $row68
$row69 = New-Object 'object[,]' 1,8
$row69[0,0] = 0
$row69[0,1] = 0
$row69[0,2] = 0
$row69[0,3] = 0
$row69[0,4] = 0
$row69[0,5] = 0
$row69[0,6] = 0
$row69[0,7] = 0
$ws.Range("D69:K69").Value2 = $row69
$row70 = New-Object 'object[,]' 1,8
$row70[0,0] = 0
$row70[0,1] = 0
$row70[0,2] = 0
$row70[0,3] = 0
$row70[0,4] = 0
$row70[0,5] = 0
$row70[0,6] = 0
$row70[0,7] = 0
$ws.Range("D70:K70").Value2 = $row70
$row71 = New-Object 'object[,]' 1,8
$row71[0,0] = 0
$row71[0,1] = 0
$row71[0,2] = 0
$row71[0,3] = 0
$row71[0,4] = 0
$row71[0,5] = 0
$row71[0,6] = 0
$row71[0,7] = 0
$ws.Range("D71:K71").Value2 = $row71
$row72 = New-Object 'object[,]' 1,8
$row72[0,0] = 34300
$row72[0,1] = 38800
$row72[0,2] = 41700
$row72[0,3] = 43400
$row72[0,4] = 46100
$row72[0,5] = 47000
$row72[0,6] = 49800
$row72[0,7] = 55000
$ws.Range("D72:K72").Value2 = $row72
$row73 = New-Object 'object[,]' 1,8
$row73[0,0] = 0
$row73[0,1] = 0
$row73[0,2] = 0
$row73[0,3] = 0
$row73[0,4] = 0
$row73[0,5] = 0
$row73[0,6] = 0
$row73[0,7] = 0
$ws.Range("D73:K73").Value2 = $row73
$row74 = New-Object 'object[,]' 1,8
$row74[0,0] = 0
$row74[0,1] = 0
$row74[0,2] = 0
$row74[0,3] = 0
$row74[0,4] = 0
$row74[0,5] = 0
$row74[0,6] = 0
$row74[0,7] = 0
$ws.Range("D74:K74").Value2 = $row74
$row75 = New-Object 'object[,]' 1,8
$row75[0,0] = 0
$row75[0,1] = 0
$row75[0,2] = 0
$row75[0,3] = 0
$row75[0,4] = 0
$row75[0,5] = 0
$row75[0,6] = 0
$row75[0,7] = 0
$ws.Range("D75:K75").Value2 = $row75
$row76 = New-Object 'object[,]' 1,8
$row76[0,0] = 39900
$row76[0,1] = 44400
$row76[0,2] = 47300
$row76[0,3] = 49200
$row76[0,4] = 51700
$row76[0,5] = 52500
$row76[0,6] = 55000
$row76[0,7] = 60400
$ws.Range("D76:K76").Value2 = $row76
$row77 = New-Object 'object[,]' 1,8
$row77[0,0] = 0
$row77[0,1] = 0
$row77[0,2] = 0
$row77[0,3] = 0
$row77[0,4] = 0
$row77[0,5] = 0
$row77[0,6] = 0
$row77[0,7] = 0
$ws.Range("D77:K77").Value2 = $row77
$row80 = New-Object 'object[,]' 1,8
$row80[0,0] = 43373
$row80[0,1] = 43281
$row80[0,2] = 43190
$row80[0,3] = 43100
$row80[0,4] = 43008
$row80[0,5] = 42916
$row80[0,6] = 42825
$row80[0,7] = 42735
$ws.Range("D80:K80").Value2 = $row80
$row81 = New-Object 'object[,]' 1,8
$row81[0,0] = -4300
$row81[0,1] = -5200
$row81[0,2] = -1800
$row81[0,3] = -2800
$row81[0,4] = -1200
$row81[0,5] = -5900
$row81[0,6] = -2700
$row81[0,7] = -3500
$ws.Range("D81:K81").Value2 = $row81
$row83 = New-Object 'object[,]' 1,8
$row83[0,0] = 1000
$row83[0,1] = 1900
$row83[0,2] = 900
$row83[0,3] = 1000
$row83[0,4] = 900
$row83[0,5] = 1600
$row83[0,6] = 800
$row83[0,7] = 800
$ws.Range("D83:K83").Value2 = $row83
$row84 = New-Object 'object[,]' 1,8
$row84[0,0] = 0
$row84[0,1] = 0
$row84[0,2] = 0
$row84[0,3] = 0
$row84[0,4] = 0
$row84[0,5] = 0
$row84[0,6] = 0
$row84[0,7] = 0
$ws.Range("D84:K84").Value2 = $row84
$row85 = New-Object 'object[,]' 1,8
$row85[0,0] = 0
$row85[0,1] = 0
$row85[0,2] = 0
$row85[0,3] = 0
$row85[0,4] = 0
$row85[0,5] = 0
$row85[0,6] = 0
$row85[0,7] = 0
$ws.Range("D85:K85").Value2 = $row85
$row86 = New-Object 'object[,]' 1,8
$row86[0,0] = 0
$row86[0,1] = 0
$row86[0,2] = 0
$row86[0,3] = 0
$row86[0,4] = 0
$row86[0,5] = 0
$row86[0,6] = 0
$row86[0,7] = 0
$ws.Range("D86:K86").Value2 = $row86
$row87 = New-Object 'object[,]' 1,8
$row87[0,0] = 0
$row87[0,1] = 0
$row87[0,2] = 0
$row87[0,3] = 0
$row87[0,4] = 0
$row87[0,5] = 0
$row87[0,6] = 0
$row87[0,7] = 0
$ws.Range("D87:K87").Value2 = $row87
$row88 = New-Object 'object[,]' 1,8
$row88[0,0] = 0
$row88[0,1] = 0
$row88[0,2] = 0
$row88[0,3] = 0
$row88[0,4] = 0
$row88[0,5] = 0
$row88[0,6] = 0
$row88[0,7] = 0
$ws.Range("D88:K88").Value2 = $row88
$row89 = New-Object 'object[,]' 1,8
$row89[0,0] = -3300
$row89[0,1] = -3800
$row89[0,2] = -200
$row89[0,3] = -600
$row89[0,4] = -3100
$row89[0,5] = -4300
$row89[0,6] = -1800
$row89[0,7] = -2200
$ws.Range("D89:K89").Value2 = $row89
$row91 = New-Object 'object[,]' 1,8
$row91[0,0] = -700
$row91[0,1] = -900
$row91[0,2] = -300
$row91[0,3] = -1400
$row91[0,4] = -100
$row91[0,5] = -2200
$row91[0,6] = -800
$row91[0,7] = 800
$ws.Range("D91:K91").Value2 = $row91
$row92 = New-Object 'object[,]' 1,8
$row92[0,0] = 0
$row92[0,1] = 0
$row92[0,2] = 0
$row92[0,3] = 0
$row92[0,4] = 0
$row92[0,5] = 0
$row92[0,6] = 0
$row92[0,7] = 0
$ws.Range("D92:K92").Value2 = $row92
$row93 = New-Object 'object[,]' 1,8
$row93[0,0] = 0
$row93[0,1] = 0
$row93[0,2] = 0
$row93[0,3] = 0
$row93[0,4] = 0
$row93[0,5] = 0
$row93[0,6] = 0
$row93[0,7] = 0
$ws.Range("D93:K93").Value2 = $row93
$row94 = New-Object 'object[,]' 1,8
$row94[0,0] = 2900
$row94[0,1] = 0
$row94[0,2] = -4400
$row94[0,3] = -4900
$row94[0,4] = 900
$row94[0,5] = -1300
$row94[0,6] = 300
$row94[0,7] = 500
$ws.Range("D94:K94").Value2 = $row94
$row96 = New-Object 'object[,]' 1,8
$row96[0,0] = 0
$row96[0,1] = 0
$row96[0,2] = 0
$row96[0,3] = 0
$row96[0,4] = 0
$row96[0,5] = 0
$row96[0,6] = 0
$row96[0,7] = 0
$ws.Range("D96:K96").Value2 = $row96
$row97 = New-Object 'object[,]' 1,8
$row97[0,0] = 0
$row97[0,1] = 0
$row97[0,2] = 0
$row97[0,3] = 0
$row97[0,4] = 0
$row97[0,5] = 0
$row97[0,6] = 0
$row97[0,7] = 0
$ws.Range("D97:K97").Value2 = $row97
$row98 = New-Object 'object[,]' 1,8
$row98[0,0] = 0
$row98[0,1] = 0
$row98[0,2] = 0
$row98[0,3] = 0
$row98[0,4] = 0
$row98[0,5] = 0
$row98[0,6] = 0
$row98[0,7] = 0
$ws.Range("D98:K98").Value2 = $row98
$row99 = New-Object 'object[,]' 1,8
$row99[0,0] = 0
$row99[0,1] = 0
$row99[0,2] = 0
$row99[0,3] = 0
$row99[0,4] = 0
$row99[0,5] = 0
$row99[0,6] = 0
$row99[0,7] = 0
$ws.Range("D99:K99").Value2 = $row99
$row100 = New-Object 'object[,]' 1,8
$row100[0,0] = -300
$row100[0,1] = -700
$row100[0,2] = -400
$row100[0,3] = 10800
$row100[0,4] = 700
$row100[0,5] = 1100
$row100[0,6] = 500
$row100[0,7] = 2400
$ws.Range("D100:K100").Value2 = $row100
$row101 = New-Object 'object[,]' 1,8
$row101[0,0] = 0
$row101[0,1] = 0
$row101[0,2] = 0
$row101[0,3] = 0
$row101[0,4] = 100
$row101[0,5] = 200
$row101[0,6] = 0
$row101[0,7] = 0
$ws.Range("D101:K101").Value2 = $row101
$row102 = New-Object 'object[,]' 1,8
$row102[0,0] = -600
$row102[0,1] = -4400
$row102[0,2] = -5000
$row102[0,3] = 5400
$row102[0,4] = -1400
$row102[0,5] = -4300
$row102[0,6] = -1000
$row102[0,7] = 800
$ws.Range("D102:K102").Value2 = $row102
